$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column G header "t" and the ply-thickness values for each material row.
$ws.Range("G1").Value = "t"
$ws.Range("G2").Value = 0.0002
$ws.Range("G3").Value = 0.0015
$ws.Range("G4").Value = 0.0003
$ws.Range("G5").Value = 0.0015

# Match the look of the existing table: G2 picks up the "Good" (green) style
# used by columns C/D, while G3:G5 pick up the "Neutral" (yellow) style used
# by column F. Copy formats from those cells instead of assigning a named
# style so the existing style records in the workbook are reused as-is.
$ws.Range("C2").Copy()
$ws.Range("G2").PasteSpecial(-4122)

$ws.Range("F3").Copy()
$ws.Range("G3:G5").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Re-assert the values in case PasteSpecial(formats) touched anything.
$ws.Range("G2").Value = 0.0002
$ws.Range("G3").Value = 0.0015
$ws.Range("G4").Value = 0.0003
$ws.Range("G5").Value = 0.0015
